# precision_summary.xlsx - fill in Binkley CNN precision/recall columns (C/D),
# add their row-14 averages, and move the active selection, per the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Binkley" precision (C) / recall (D) data, rows 4-13 ---
# (ordered array of row/C/D triples so write order is deterministic)
$values = @(
    , @(4,  0.79318239999999995, 0.79071619999999998)
    , @(5,  0.78846156999999994, 0.7820513)
    , @(6,  0.79707205000000003, 0.80780779999999996)
    , @(7,  0.76423909999999995, 0.77964520000000004)
    , @(8,  0.79358620000000002, 0.79900629999999995)
    , @(9,  0.79739444999999998, 0.78032345000000003)
    , @(10, 0.80216799999999999, 0.79403800000000002)
    , @(11, 0.77959920000000005, 0.73330289999999998)
    , @(12, 0.82234436,          0.80540290000000003)
    , @(13, 0.78838956000000004, 0.81039333000000002)
)

foreach ($triple in $values) {
    $row = $triple[0]
    $ws.Cells.Item($row, 3).Value = $triple[1]
    $ws.Cells.Item($row, 4).Value = $triple[2]
}

# --- Row 14 averages for the newly-filled columns, matching the existing
#     E14/F14 "Average" row look (plain Times New Roman, no border) ---
$ws.Range("E14").Copy() | Out-Null
$ws.Range("C14:D14").PasteSpecial(-4122) | Out-Null
$ws.Range("C14:D14").Formula = "=AVERAGE(C4:C13)"

# --- Move the active cell / selection ---
$ws.Range("J16").Select() | Out-Null
